# Update "想去人数" (interest count) figures refreshed by the data scrape,
# per commit "Update gh-pages to output generated at 74db155".
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 57
$wsExhibit.Range("F3").Value = 979
$wsExhibit.Range("F5").Value = 10876
$wsExhibit.Range("F7").Value = 362
$wsExhibit.Range("F9").Value = 1947
$wsExhibit.Range("F12").Value = 228
$wsExhibit.Range("F16").Value = 990
$wsExhibit.Range("F20").Value = 422

# Sheet "演出" (performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F7").Value = 633

# Sheet "全部类型" (all types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 57
$wsAll.Range("F5").Value = 979
$wsAll.Range("F8").Value = 10876
$wsAll.Range("F11").Value = 362
$wsAll.Range("F13").Value = 1947
$wsAll.Range("F17").Value = 228
$wsAll.Range("F21").Value = 990
$wsAll.Range("F23").Value = 633
$wsAll.Range("F26").Value = 422
